$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 297.5
$ws.Range("I2").Value = 267.85715
$ws.Range("J2").Value = 366.66666
$ws.Range("K2").Value = 267.85715
$ws.Range("L2").Value = 366.66666
$ws.Range("M2").Value = -154.85715
$ws.Range("N2").Value = -592.66666

$ws.Range("H28").Value = 502.9375
$ws.Range("I28").Value = 524.5172
$ws.Range("J28").Value = 294.33334
$ws.Range("K28").Value = 524.5172
$ws.Range("L28").Value = 294.33334
$ws.Range("M28").Value = -39.5172
$ws.Range("N28").Value = -1264.33334

$ws.Range("H64").Value = 5332.6665
$ws.Range("J64").Value = 5999
$ws.Range("L64").Value = 5999
$ws.Range("N64").Value = -6495

$ws.Range("H67").Value = 5332.6665
$ws.Range("J67").Value = 5999
$ws.Range("L67").Value = 5999
$ws.Range("N67").Value = -7715

$ws.Range("H70").Value = 1333
$ws.Range("J70").Value = 999.5
$ws.Range("L70").Value = 2998.5
$ws.Range("N70").Value = -3538.5

$ws.Range("H73").Value = 1333
$ws.Range("J73").Value = 999.5
$ws.Range("L73").Value = 2998.5
$ws.Range("N73").Value = -4870.5

$ws.Range("H112").Value = 9463.08
$ws.Range("J112").Value = 8609.315
$ws.Range("L112").Value = 25827.945
$ws.Range("N112").Value = -28043.945

$ws.Range("H116").Value = 4664.222
$ws.Range("J116").Value = 3759.625
$ws.Range("L116").Value = 3759.625
$ws.Range("N116").Value = -10643.625

$ws.Range("H132").Value = 9631.037
$ws.Range("I132").Value = 6566.522
$ws.Range("K132").Value = 19699.566
$ws.Range("M132").Value = -17169.566

$ws.Range("H135").Value = 1139.4286
$ws.Range("I135").Value = 1162.8334
$ws.Range("K135").Value = 10465.5006
$ws.Range("M135").Value = -7930.500599999999

$ws.Range("H137").Value = 15127
$ws.Range("J137").Value = 20017.53
$ws.Range("L137").Value = 60052.59
$ws.Range("N137").Value = -65152.59

$ws.Range("H138").Value = 4024.6353
$ws.Range("J138").Value = 4084.4736
$ws.Range("L138").Value = 12253.4208
$ws.Range("N138").Value = -22533.4208

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()

$ws.Range("H25").Value = 593.6
$ws.Range("I25").Value = 593.6
$ws.Range("K25").Value = 593.6
$ws.Range("M25").Value = -191.6

$ws.Range("H32").Value = 10740.169
$ws.Range("I32").Value = 3138.7805
$ws.Range("K32").Value = 3138.7805
$ws.Range("M32").Value = -2851.7805

$ws.Range("H35").Value = 4605
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()

$ws.Range("H61").Value = 13206.363
$ws.Range("I61").Value = 7420.8
$ws.Range("J61").Value = 22107.23
$ws.Range("K61").Value = 7420.8
$ws.Range("L61").Value = 22107.23
$ws.Range("M61").Value = -7208.8
$ws.Range("N61").Value = -22531.23

$ws.Range("H63").Value = 2596
$ws.Range("I63").Value = 2596
$ws.Range("K63").Value = 2596
$ws.Range("M63").Value = -1910

$ws.Range("H66").Value = 2596
$ws.Range("I66").Value = 2596
$ws.Range("K66").Value = 12980
$ws.Range("M66").Value = -9548

$ws.Range("H97").Value = 1391.6487
$ws.Range("I97").Value = 1525.9375
$ws.Range("K97").Value = 1525.9375
$ws.Range("M97").Value = -1029.9375

$ws.Range("H110").Value = 9914.5
$ws.Range("I110").Value = 11097.4
$ws.Range("K110").Value = 11097.4
$ws.Range("M110").Value = -9052.4

$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()

$ws.Range("H132").Value = 11079815
$ws.Range("I132").Value = 57096.715
$ws.Range("K132").Value = 171290.145
$ws.Range("M132").Value = -168760.145

$ws.Range("H135").Value = 132404.83
$ws.Range("J135").Value = 132404.83
$ws.Range("L135").Value = 132404.83
$ws.Range("N135").Value = -142544.83

$ws.Range("H136").Value = 13206.363
$ws.Range("I136").Value = 7420.8
$ws.Range("J136").Value = 22107.23
$ws.Range("K136").Value = 22262.4
$ws.Range("L136").Value = 66321.69
$ws.Range("M136").Value = -19712.4
$ws.Range("N136").Value = -71421.69

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value = 4196.2856
$ws.Range("I54").Value = 3875.8
$ws.Range("J54").Value = 4997.5
$ws.Range("K54").Value = 3875.8
$ws.Range("L54").Value = 4997.5
$ws.Range("M54").Value = -3391.8
$ws.Range("N54").Value = -5965.5

$ws.Range("H81").Value = 31279.834
$ws.Range("I81").Value = 14995
$ws.Range("J81").Value = 34536.8
$ws.Range("K81").Value = 14995
$ws.Range("L81").Value = 34536.8
$ws.Range("M81").Value = -13934
$ws.Range("N81").Value = -36658.8

$ws.Range("H84").Value = 31279.834
$ws.Range("I84").Value = 14995
$ws.Range("J84").Value = 34536.8
$ws.Range("K84").Value = 44985
$ws.Range("L84").Value = 103610.4
$ws.Range("M84").Value = -39681
$ws.Range("N84").Value = -114218.4

$ws.Range("H94").Value = 1955.75
$ws.Range("I94").Value = 1582.7693
$ws.Range("J94").Value = 2648.4285
$ws.Range("K94").Value = 1582.7693
$ws.Range("L94").Value = 2648.4285
$ws.Range("M94").Value = -1131.7693
$ws.Range("N94").Value = -3550.4285

$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("M97").ClearContents()

$ws.Range("H102").Value = 11755.818
$ws.Range("I102").Value = 11724.7
$ws.Range("J102").Value = 12067
$ws.Range("K102").Value = 11724.7
$ws.Range("L102").Value = 12067
$ws.Range("M102").Value = -8479.7
$ws.Range("N102").Value = -18557

$ws.Range("H134").Value = 35017.332
$ws.Range("I134").Value = 22522.5
$ws.Range("K134").Value = 67567.5
$ws.Range("M134").Value = -65032.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 1663.25
$ws.Range("I5").Value = 301.55554
$ws.Range("J5").Value = 5748.3335
$ws.Range("K5").Value = 301.55554
$ws.Range("L5").Value = 5748.3335
$ws.Range("M5").Value = -189.55554
$ws.Range("N5").Value = -5972.3335

$ws.Range("H12").Value = 1402
$ws.Range("I12").Value = 1482.4
$ws.Range("K12").Value = 1482.4
$ws.Range("M12").Value = -1312.4

$ws.Range("H31").Value = 69808.375
$ws.Range("I31").Value = 185266.1
$ws.Range("J31").Value = 20960.885
$ws.Range("K31").Value = 185266.1
$ws.Range("L31").Value = 20960.885
$ws.Range("M31").Value = -184971.1
$ws.Range("N31").Value = -21550.885

$ws.Range("H34").Value = 69808.375
$ws.Range("I34").Value = 185266.1
$ws.Range("J34").Value = 20960.885
$ws.Range("K34").Value = 185266.1
$ws.Range("L34").Value = 20960.885
$ws.Range("M34").Value = -185064.1
$ws.Range("N34").Value = -21364.885

$ws.Range("H105").Value = 21602
$ws.Range("I105").Value = 100010
$ws.Range("J105").Value = 2000
$ws.Range("K105").Value = 100010
$ws.Range("L105").Value = 2000
$ws.Range("M105").Value = -98263
$ws.Range("N105").Value = -5494

$ws.Range("H117").Value = 250000
$ws.Range("I117").Value = 0
$ws.Range("K117").Value = 0
$ws.Range("M117").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 13402605
$ws.Range("J32").Value = 11113976
$ws.Range("L32").Value = 33341928
$ws.Range("N32").Value = -33342494

$ws.Range("H68").Value = 1834.44
$ws.Range("J68").Value = 1924.6086
$ws.Range("L68").Value = 5773.825800000001
$ws.Range("N68").Value = -7395.825800000001

$ws.Range("H71").Value = 1834.44
$ws.Range("J71").Value = 1924.6086
$ws.Range("L71").Value = 17321.4774
$ws.Range("N71").Value = -25433.4774

$ws.Range("H116").Value = 12650596
$ws.Range("I116").Value = 16866378
$ws.Range("K116").Value = 50599134
$ws.Range("M116").Value = -50595692

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 18333.334
$ws.Range("J21").Value = 18333.334
$ws.Range("L21").Value = 18333.334
$ws.Range("N21").Value = -18679.334

$ws.Range("H30").Value = 18333.334
$ws.Range("J30").Value = 18333.334
$ws.Range("L30").Value = 18333.334
$ws.Range("N30").Value = -18543.334

$ws.Range("H31").Value = 749.5
$ws.Range("I31").Value = 749.5
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 749.5
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -457.5
$ws.Range("N31").ClearContents()

$ws.Range("H37").Value = 749.5
$ws.Range("I37").Value = 749.5
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 749.5
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -472.5
$ws.Range("N37").ClearContents()

$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()

$ws.Range("H55").Value = 5422.5
$ws.Range("J55").Value = 5500
$ws.Range("L55").Value = 5500
$ws.Range("N55").Value = -6154

$ws.Range("H102").Value = 6140.684
$ws.Range("I102").Value = 6891.1875
$ws.Range("J102").Value = 2138
$ws.Range("K102").Value = 6891.1875
$ws.Range("L102").Value = 2138
$ws.Range("M102").Value = -5269.1875
$ws.Range("N102").Value = -5382

$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()

$ws.Range("H122").Value = 2705.45
$ws.Range("I122").Value = 2340.7334
$ws.Range("J122").Value = 3799.6
$ws.Range("K122").Value = 7022.2002
$ws.Range("L122").Value = 11398.8
$ws.Range("M122").Value = -4572.2002
$ws.Range("N122").Value = -16298.8

$ws.Range("H126").Value = 6475.2
$ws.Range("I126").Value = 7162.6924
$ws.Range("J126").Value = 2006.5
$ws.Range("K126").Value = 21488.0772
$ws.Range("L126").Value = 6019.5
$ws.Range("M126").Value = -19018.0772
$ws.Range("N126").Value = -10959.5

$ws.Range("H132").Value = 289313.38
$ws.Range("I132").Value = 5457.841
$ws.Range("K132").Value = 16373.523
$ws.Range("M132").Value = -13843.523

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 16447.9
$ws.Range("J7").Value = 8000
$ws.Range("L7").Value = 8000
$ws.Range("N7").Value = -8224

$ws.Range("H68").Value = 15104.368
$ws.Range("I68").Value = 17385.215
$ws.Range("J68").Value = 8718
$ws.Range("K68").Value = 17385.215
$ws.Range("L68").Value = 8718
$ws.Range("M68").Value = -16636.215
$ws.Range("N68").Value = -10216

$ws.Range("H71").Value = 15104.368
$ws.Range("I71").Value = 17385.215
$ws.Range("J71").Value = 8718
$ws.Range("K71").Value = 86926.075
$ws.Range("L71").Value = 43590
$ws.Range("M71").Value = -83182.075
$ws.Range("N71").Value = -51078

$ws.Range("H93").Value = 9908.167
$ws.Range("I93").Value = 10285.177
$ws.Range("J93").Value = 3499
$ws.Range("K93").Value = 10285.177
$ws.Range("L93").Value = 3499
$ws.Range("M93").Value = -9037.177
$ws.Range("N93").Value = -5995

$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()

$ws.Range("H126").Value = 16447.9
$ws.Range("J126").Value = 8000
$ws.Range("L126").Value = 24000
$ws.Range("N126").Value = -28940

$ws.Range("H136").Value = 119045.25
$ws.Range("I136").Value = 24559.334
$ws.Range("J136").Value = 196351.9
$ws.Range("K136").Value = 73678.002
$ws.Range("L136").Value = 589055.7
$ws.Range("M136").Value = -71128.002
$ws.Range("N136").Value = -594155.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 6499.25

$ws.Range("H62").Value = 33499.75
$ws.Range("I62").Value = 7497
$ws.Range("J62").Value = 42167.332
$ws.Range("K62").Value = 7497
$ws.Range("L62").Value = 42167.332
$ws.Range("M62").Value = -6873
$ws.Range("N62").Value = -43415.332

$ws.Range("H65").Value = 33499.75
$ws.Range("I65").Value = 7497
$ws.Range("J65").Value = 42167.332
$ws.Range("K65").Value = 37485
$ws.Range("L65").Value = 210836.66
$ws.Range("M65").Value = -34365
$ws.Range("N65").Value = -217076.66

$ws.Range("H103").Value = 0
$ws.Range("I103").Value = 0
$ws.Range("K103").Value = 0
$ws.Range("M103").ClearContents()

$ws.Range("H109").Value = 23000
$ws.Range("J109").Value = 23000
$ws.Range("L109").Value = 23000
$ws.Range("N109").Value = -25774

$ws.Range("H113").Value = 3191.0527
$ws.Range("I113").Value = 3818.8
$ws.Range("J113").Value = 837
$ws.Range("K113").Value = 11456.4
$ws.Range("L113").Value = 2511
$ws.Range("M113").Value = -9286.400000000001
$ws.Range("N113").Value = -6851

$ws.Range("H122").Value = 4174.049
$ws.Range("I122").Value = 3657.5625
$ws.Range("J122").Value = 6010.4443
$ws.Range("K122").Value = 10972.6875
$ws.Range("L122").Value = 18031.3329
$ws.Range("M122").Value = -8522.6875
$ws.Range("N122").Value = -22931.3329

$ws.Range("H136").Value = 29862.105
$ws.Range("I136").Value = 3113.4167
$ws.Range("J136").Value = 75717
$ws.Range("K136").Value = 9340.250100000001
$ws.Range("L136").Value = 227151
$ws.Range("M136").Value = -6790.250100000001
$ws.Range("N136").Value = -232251
